$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 112 timestamps (F112, G112) ---
$ws.Range("F112").Value = 45910.65059689815
$ws.Range("G112").Value = 45910.65059133102

# --- Data for newly appended rows 113-128 ---
$newRows = @(
    ,@{A='Fucntionality_test_Letter-high_with_(5)-NN_Classifier_GED'; B='Letter-high'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_Letter-high.joblib'; F=45912.60330827547; G=45912.60326099537}
    ,@{A='Fucntionality_test_Letter-high_with_(5)-NN_Classifier_GED'; B='Letter-high'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_Letter-high.joblib'; F=45912.6047197338; G=45912.60468278935}
    ,@{A='Fucntionality_test_Letter-high_with_(5)-NN_Classifier_GED'; B='Letter-high'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_Letter-high.joblib'; F=45912.60630778935; G=45912.60627363426}
    ,@{A='Fucntionality_test_Letter-high_with_(5)-NN_Classifier_GED'; B='Letter-high'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_Letter-high.joblib'; F=45912.60910956019; G=45912.60907547454}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.61058741898; G=45912.6105871875}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.61173458333; G=45912.61173435185}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.61266902777; G=45912.61266877314}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.612916875; G=45912.61291667824}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.612916875; G=45912.61291667824}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.61332822917; G=45912.61332799769}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.61519488426; G=45912.61519466435}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.61885483797; G=45912.61885425926}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.65649662037; G=45912.65649607639}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.65649662037; G=45912.65649607639}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.65889356482; G=45912.65889335648}
    ,@{A='Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED'; B='MUTAG'; C='(5)-NN_Classifier_GED'; D='(5)-NN_Classifier_GED_trained_on_MUTAG.joblib'; F=45912.65889356853; G=45912.65889335222}
)

$startRow = 113
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    # Columns E (dataset_load_duration) and H (Error) stay empty for these rows.
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G

    # Match the number format used by the existing timestamp columns (F, G)
    $ws.Cells.Item($r, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
